# Applies the Mura_List_Template.xlsx edits:
#  1. Rename header cells C1/D1 from "Before/After Adj (87.5IRE)" to "Before/After Adj"
#  2. Fill in the continuing "No" sequence in column A for rows 23..95 (values 22..94)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mura")

# 1. Update header text in C1 / D1
$ws.Range("C1").Value = "Before Adj"
$ws.Range("D1").Value = "After Adj"

# 2. Number column A continuing from row 22 (value 21) down through row 95 (value 94)
for ($row = 23; $row -le 95; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 1
}
